$d = $word.ActiveDocument

$replacements = @(
    @("2025-07-02 Wednesday", "2025-07-03 Thursday"),
    @("167×6=", "809×9="),
    @("155×8=", "778×4="),
    @("549×9=", "949×4="),
    @("978×3=", "111×8="),
    @("783×9=", "109×9="),
    @("374×8=", "877×5="),
    @("703×2=", "114×7="),
    @("556×7=", "315×9="),
    @("648×4=", "982×9="),
    @("638×9=", "552×6="),
    @("207×6=", "786×2="),
    @("978×7=", "396×8="),
    @("143×2=", "442×9="),
    @("392×8=", "456×4="),
    @("605×9=", "365×4="),
    @("749×8=", "855×3="),
    @("548×4=", "992×2="),
    @("731×3=", "619×7="),
    @("782×9=", "955×2="),
    @("607×8=", "819×6="),
    @("682×7=", "912×7="),
    @("635×9=", "842×3="),
    @("378×5=", "925×4="),
    @("298×8=", "313×6="),
    @("720×4=", "908×8=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
